$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A3: "ID 01" -> "ID 02" (adds a new shared string)
$ws.Range("A3").Value = "ID 02"

# Update C2: 5 -> 50
$ws.Range("C2").Value = 50

# Update C3: 15 -> 50
$ws.Range("C3").Value = 50
